# Applies the coinranking-data refresh described in the commit:
#   "Updated cryptos list on Mon Oct 28 13:30:43 UTC 2024 with GitHub Actions"
#
# Price/volume figures are refreshed in place, and a handful of coins swap
# table rows (their rank position stayed the same, the coin occupying it
# changed), which also means their Coin/Link cells must be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below hold numeric-looking figures that must stay literal text
# (e.g. "0.520", "9.00") -- without forcing a Text number format first,
# Excel would parse them as numbers and silently drop the trailing zero.
$textCells = @("D5", "D6", "D8", "D10", "D14", "D17", "D19", "D21", "D22", "D24", "D27", "D29", "D30", "D32", "D34", "D37", "D38", "D39", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}


$ws.Range("D2").Value = "69.005.09"
$ws.Range("E2").Value = "  +1.93%  "

$ws.Range("D3").Value = "2.529.80"
$ws.Range("E3").Value = "  +1.24%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "596.05"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").Value = "176.85"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("D9").Value = "2.529.54"
$ws.Range("E9").Value = "  +1.17%  "

$ws.Range("D10").Value = "0.148"
$ws.Range("E10").Value = "  +5.93%  "

$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("D14").Value = "26.24"
$ws.Range("E14").Value = "  +2.00%  "

$ws.Range("D15").Value = "2.951.35"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").Value = "68.898.32"
$ws.Range("E16").Value = "  +2.04%  "

$ws.Range("D17").Value = "0.0000174"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").Value = "2.527.14"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("D19").Value = "11.15"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").Value = "361.03"
$ws.Range("E21").Value = "  +2.45%  "

$ws.Range("D22").Value = "4.18"
$ws.Range("E22").Value = "  +2.81%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "70.81"
$ws.Range("E24").Value = "  +0.45%  "

$ws.Range("E25").Value = "  -0.52%  "

$ws.Range("E26").Value = "  -6.45%  "

$ws.Range("D27").Value = "9.00"
$ws.Range("E27").Value = "  -3.15%  "

$ws.Range("D28").Value = "2.654.78"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").Value = "521.53"
$ws.Range("E30").Value = "  +2.67%  "

$ws.Range("D31").Value = "0.0₃0895"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("D32").Value = "7.78"
$ws.Range("E32").Value = "  -1.03%  "

$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  +0.68%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("D37").Value = "162.80"
$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "18.52"
$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "18.69"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("E40").Value = "  +4.02%  "

$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "0.328"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").Value = "4.85"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("D45").Value = "2.40"
$ws.Range("E45").Value = "  -1.23%  "

$ws.Range("D46").Value = "151.92"
$ws.Range("E46").Value = "  +5.48%  "

$ws.Range("D47").Value = "3.60"
$ws.Range("E47").Value = "  +2.72%  "

$ws.Range("E48").Value = "  +1.12%  "

$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "1.61"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0741"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.581"
$ws.Range("E51").Value = "  -1.15%  "
